# Applies the "reading excel and applying discount" edit:
#  - Updates the client code footnote with the client's actual email address
#  - Replaces the (placeholder) discount unit-price cell with a note stating
#    that this client does not benefit from any discount, which in turn makes
#    the dependent formulas (line total, subtotal, tax, total) evaluate as
#    errors because they now try to multiply by text instead of a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# CLIENT CODE note (was "name, email address")
$ws.Range("A31").Value = "charlie.charlie@mail.com"

# Client discount line - unit price cell now holds an explanatory note
# instead of a numeric value.
$ws.Range("E18").Value = "This client doesn't benefit from any discount"

$excel.Calculate()
